# Update "想去人数" (want-to-go count) values in column F across sheets,
# matching the upstream data refresh recorded in the commit.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 499
$ws1.Range("F3").Value = 1520
$ws1.Range("F4").Value = 796
$ws1.Range("F5").Value = 205
$ws1.Range("F8").Value = 693
$ws1.Range("F10").Value = 1349
$ws1.Range("F12").Value = 1012
$ws1.Range("F13").Value = 20
$ws1.Range("F14").Value = 56
$ws1.Range("F16").Value = 42
$ws1.Range("F20").Value = 530
$ws1.Range("F21").Value = 549
$ws1.Range("F24").Value = 164

# 演出 (Performances) sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 237
$ws2.Range("F10").Value = 56

# 本地生活 (Local life) sheet
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 182

# 全部类型 (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 499
$ws4.Range("F3").Value = 182
$ws4.Range("F4").Value = 1520
$ws4.Range("F6").Value = 796
$ws4.Range("F7").Value = 205
$ws4.Range("F11").Value = 693
$ws4.Range("F13").Value = 1349
$ws4.Range("F15").Value = 1012
$ws4.Range("F16").Value = 20
$ws4.Range("F17").Value = 56
$ws4.Range("F19").Value = 42
$ws4.Range("F22").Value = 237
$ws4.Range("F27").Value = 530
$ws4.Range("F28").Value = 549
$ws4.Range("F32").Value = 164
$ws4.Range("F34").Value = 56
$ws4.Range("F35").Value = 56
